$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Weekly Timesheet")

# Row 2 - Jones (was Campbell)
$ws.Range("B2").Value = "Jones"
$ws.Range("C2").Value = 9.5
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 950

# Row 3 - Lynn (was Regan)
$ws.Range("B3").Value = "Lynn"
$ws.Range("C3").Value = 7.5
$ws.Range("E3").Value = 90
$ws.Range("F3").Value = 675

# Row 4 - Davis (was Fritts); numeric values unchanged
$ws.Range("B4").Value = "Davis"

# Row 5 - Goodrich (was Zygmunt)
$ws.Range("B5").Value = "Goodrich"
$ws.Range("C5").Value = 9
$ws.Range("F5").Value = 900

# Row 6 - "Total Hours on Insp." (was Keevil)
$ws.Range("B6").Value = "Total Hours on Insp."
$ws.Range("C6").Value = 8
$ws.Range("F6").Value = 800

# Row 8 - SUBTOTAL row
$ws.Range("C8").Value = 42.5
$ws.Range("D8").Value = "Reg: 42.5 / OT: 0"
$ws.Range("F8").Value = 4175
